{"js": "// The edit:\n//  1. Removes the old \"_GoBack\" bookmark that sat at the end of the\n//     \"My mine belong belongs his --- people countries commonwealth\" paragraph.\n//  2. Replaces the second of the two blank paragraphs that precede\n//     \"Partly to do with first-decade referring to people as \"my people\":\"\n//     with three new paragraphs of notes, the first of which ends with a\n//     (new) \"_GoBack\" bookmark.\n\n// Step 1: drop the stray _GoBack bookmark left over from the previous\n// cursor position. Office.js surfaces this the same way Word's COM object\n// model does: Document.deleteBookmark(name).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Step 2: locate the two consecutive empty paragraphs that sit right\n// before the \"Partly to do with first-decade...\" paragraph; the *second*\n// one is the insertion point that gets replaced by the new content.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length - 1; i++) {\n  const cur = paragraphs.items[i];\n  const next = paragraphs.items[i + 1];\n  if (\n    cur.text === \"\" &&\n    next.text.indexOf(\"Partly to do with first-decade referring to people as\") === 0\n  ) {\n    target = cur;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the target blank paragraph.\");\n}\n\n// Build the replacement OOXML (three paragraphs) and insert it in place of\n// the blank paragraph via a collapsed range at its start. A package-style\n// insertOoxml call preserves the exact run / bookmark / proofErr structure\n// (unlike insertParagraph/insertText, which normalize adjacent runs that\n// share identical formatting).\nconst newParagraphsXml =\n  \"<w:p>\" +\n  \"<w:r><w:t>Investigate into this\\u2026</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> what causes greater use of my /  mine early on?</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n  \"</w:p>\" +\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">Looking at what co-occurs with </w:t></w:r>' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  \"<w:r><w:t>my</w:t></w:r>\" +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  \"</w:p>\" +\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">Earlier on </w:t></w:r>' +\n  '<w:r><w:rPr><w:u w:val=\"single\"/></w:rPr><w:t>is my people</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  \"</w:p>\";\n\nconst ooxmlPackage =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  newParagraphsXml +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nconst insertionPoint = target.getRange(\"Start\");\ninsertionPoint.insertOoxml(ooxmlPackage, \"Before\");\nawait context.sync();\n", "ps1": "# The edit:\n#  1. Removes the old \"_GoBack\" bookmark that sat at the end of the\n#     \"My mine belong belongs his --- people countries commonwealth\" paragraph.\n#  2. Replaces the second of the two blank paragraphs that precede\n#     \"Partly to do with first-decade referring to people as \"my people\":\"\n#     with three new paragraphs of notes, the first of which ends with a\n#     (new) \"_GoBack\" bookmark.\n\n$d = $word.ActiveDocument\n\n# Step 1: drop the stray _GoBack bookmark left over from the previous\n# cursor position.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Step 2: locate the two consecutive empty paragraphs that sit right\n# before the \"Partly to do with first-decade...\" paragraph; the *second*\n# one is the insertion point that gets replaced by the new content.\n$target = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le ($count - 1); $i++) {\n    $cur = $d.Paragraphs.Item($i)\n    $next = $d.Paragraphs.Item($i + 1)\n    if ($cur.Range.Text -eq \"`r\" -and $next.Range.Text.StartsWith(\"Partly to do with first-decade referring to people as\")) {\n        $target = $cur\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not locate the target blank paragraph.\"\n}\n\n# Build the replacement WordprocessingML (three paragraphs) and insert it in\n# place of the blank paragraph via Range.InsertXML. This preserves the exact\n# run / bookmark / proofErr structure (unlike Range.InsertAfter, which would\n# normalize everything into a single plain run).\n$ns = \"xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'\"\n\n$newParagraphsXml =\n    \"<w:p $ns>\" +\n    \"<w:r><w:t>Investigate into this&#8230;</w:t></w:r>\" +\n    \"<w:r><w:t xml:space='preserve'> what causes greater use of my /  mine early on?</w:t></w:r>\" +\n    \"<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>\" +\n    \"</w:p>\" +\n    \"<w:p $ns>\" +\n    \"<w:r><w:t xml:space='preserve'>Looking at what co-occurs with </w:t></w:r>\" +\n    \"<w:proofErr w:type='gramStart'/>\" +\n    \"<w:r><w:t>my</w:t></w:r>\" +\n    \"<w:proofErr w:type='gramEnd'/>\" +\n    \"</w:p>\" +\n    \"<w:p $ns>\" +\n    \"<w:r><w:t xml:space='preserve'>Earlier on </w:t></w:r>\" +\n    \"<w:r><w:rPr><w:u w:val='single'/></w:rPr><w:t>is my people</w:t></w:r>\" +\n    \"<w:r><w:t xml:space='preserve'> </w:t></w:r>\" +\n    \"</w:p>\"\n\n$target.Range.InsertXML($newParagraphsXml)\n"}
